$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: Categoría -> Categoría Y Tipo
$ws.Range("B1").Value = "Categoría Y Tipo"

# Update row 2 category value
$ws.Range("B2").Value = "pequeños Automóvil"

# Update row 3 category value
$ws.Range("B3").Value = "SUV Automóvil"

# Update row 7 category value (was "pequeños")
$ws.Range("B7").Value = "pequeños Automóvil"

# Add new row 8 with a new rental record
$ws.Range("A8").Value = "7"
$ws.Range("B8").Value = "Urbana Bicicleta"
$ws.Range("C8").Value = "15/12/2023"
$ws.Range("D8").Value = "Ruedas"
$ws.Range("E8").Value = "Exostos"
$ws.Range("F8").Value = "16/12/2023"
$ws.Range("G8").Value = "Juan"
$ws.Range("H8").Value = "Gen"
$ws.Range("I8").Value = "0"
$ws.Range("J8").Value = "No"
$ws.Range("K8").Value = "No"
$ws.Range("L8").Value = "No"
$ws.Range("M8").Value = "3000"
